# V15 - Alchemist sterker.
# The Alchemist no longer removes himself; it now captures an adjacent
# piece without moving. Update both the English and Dutch reference rows,
# and move the stray "_GoBack" bookmark so it again sits right after the
# (now edited) English description, as produced by a real Word editing
# session.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. English table: "Can remove both himself and an adjacent piece from
#    the board." -> "Can capture any adjacent piece without moving."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Can remove both himself and an adjacent piece from the board.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Can capture any adjacent piece without moving.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Re-anchor the "_GoBack" bookmark immediately after the run we just
#    edited (that's where Word leaves it after the last edit). A
#    zero-length range sitting exactly at a paragraph's end is rejected
#    by Bookmarks.Add, so nudge the boundary out with a throw-away
#    character, drop the bookmark, then remove the character again.
# ---------------------------------------------------------------------
$editedRun = $d.Content
$editedRun.Find.Execute("Can capture any adjacent piece without moving.") | Out-Null
$editedRun.Collapse(0)
$editedRun.InsertAfter("X")
$bookmarkSpot = $d.Range($editedRun.Start, $editedRun.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)
$scratchChar = $d.Range($editedRun.Start, $editedRun.Start + 1)
$scratchChar.Delete()

# ---------------------------------------------------------------------
# 3. Dutch table: split "Kan zowel zichzelf als een aanliggend stuk
#    verwijderen." into three runs reading "Kan elk" + " aanliggend
#    stuk" + " slaan zonder te bewegen."
# ---------------------------------------------------------------------
$nlPart1 = "Kan elk"
$nlPart2 = " aanliggend stuk"
$nlPart3 = " slaan zonder te bewegen."

$nlRange = $d.Content
$nlRange.Find.Execute("Kan zowel zichzelf als een aanliggend stuk verwijderen.") | Out-Null
$nlRange.Text = $nlPart1 + $nlPart2 + $nlPart3
$nlStart = $nlRange.Start
$nlEnd = $nlRange.End

# Force a run split at each boundary by flipping a character attribute
# off then back on; same-value re-writes are no-ops in this runtime, but
# toggling genuinely splits the surrounding text into its own run.
# Split the trailing " slaan zonder te bewegen." off first (working
# right-to-left keeps the earlier offsets valid).
$nlTail = $d.Range($nlStart + $nlPart1.Length + $nlPart2.Length, $nlEnd)
$nlTail.Font.Italic = $false
$nlTail.Font.Italic = $true

# Then split " aanliggend stuk" off from the remaining "Kan elk" run.
$nlMiddle = $d.Range($nlStart + $nlPart1.Length, $nlStart + $nlPart1.Length + $nlPart2.Length)
$nlMiddle.Font.Italic = $false
$nlMiddle.Font.Italic = $true
